$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73; this shifts existing rows 73-107 down to 74-108
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with its data.
# Columns A-C, E-I, N-O, Q-R keep the same values the "template" row used to have,
# while D, J, K, L, M, P contain the new data from the commit.
$ws.Range("A73").Value = 4
$ws.Range("B73").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C73").Value = "Los Lagos"
$ws.Range("D73").Value = 44582
$ws.Range("E73").Value = 10
$ws.Range("F73").Value = 100112022
$ws.Range("G73").Value = "Arveja Verde"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 70
$ws.Range("K73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = 30000
$ws.Range("N73").Value = "$/saco 25 kilos"
$ws.Range("O73").Value = "Región de La Araucanía"
$ws.Range("P73").Value = 1200
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the other date cells in column D.
$ws.Range("D73").NumberFormat = $ws.Range("D74").NumberFormat()
